$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8 gets its own (non-shared) IF formula.
$ws.Range("C8").Formula = '=IF(B8="N",A8*$B$4,A8*$B$3)+$B$2'

# Rows 9-13 are filled as a single shared-formula block (e.g. via fill-down).
$ws.Range("C9:C13").Formula = '=IF(B9="N",A9*$B$4,A9*$B$3)+$B$2'

# A stray currency-formatted cell further down/right on the sheet.
$ws.Range("F15").NumberFormat = '_-* #,##0.00\ "€"_-;\-* #,##0.00\ "€"_-;_-* "-"??\ "€"_-;_-@_-'

# Final selection ends up on C8.
[void]$ws.Range("C8").Select()
